$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that receive a new "checked" mark (✓) in column K, reusing the
# existing "好" (Good) cell style already present in the workbook.
$rows = @(2,4,5,6,7,8,10,11,12,13,14,15,18,23,24,27,29,33,36,38,41,43,46,47)

foreach ($r in $rows) {
    $cell = $ws.Range("K$r")
    $cell.Value = [char]0x2713
    $cell.Style = "好"
}

# Update the saved view state: selecting K48 both clears the stale
# topLeftCell="A7" scroll anchor and moves the active cell/selection to K48.
$ws.Range("K48").Select()
